# Scheduled-runner style refresh of market-price-derived columns
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# across the per-job leve profit tables. Only H..N value cells move;
# leve metadata (A..G) is untouched.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 51041.95
$ws.Range("J17").Value = 53691.527
$ws.Range("L17").Value = 161074.581
$ws.Range("N17").Value = -161410.581
# Row 64
$ws.Range("H64").Value = 3269.7
$ws.Range("I64").Value = 2924.25
$ws.Range("K64").Value = 2924.25
$ws.Range("M64").Value = -2676.25
# Row 67
$ws.Range("H67").Value = 3269.7
$ws.Range("I67").Value = 2924.25
$ws.Range("K67").Value = 2924.25
$ws.Range("M67").Value = -2066.25
# Row 69
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
# Row 72
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
# Row 76
$ws.Range("H76").Value = 11500646
$ws.Range("I76").Value = 7564.4546
$ws.Range("J76").Value = 47621756
$ws.Range("K76").Value = 7564.4546
$ws.Range("L76").Value = 47621756
$ws.Range("M76").Value = -7249.4546
$ws.Range("N76").Value = -47622386
# Row 79
$ws.Range("H79").Value = 11500646
$ws.Range("I79").Value = 7564.4546
$ws.Range("J79").Value = 47621756
$ws.Range("K79").Value = 7564.4546
$ws.Range("L79").Value = 47621756
$ws.Range("M79").Value = -6472.4546
$ws.Range("N79").Value = -47623940
# Row 98
$ws.Range("H98").Value = 156250540
$ws.Range("I98").Value = 178571620
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 178571620
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = -178570122
$ws.Range("N98").Value = -5996
# Row 122
$ws.Range("H122").Value = 156250540
$ws.Range("I122").Value = 178571620
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 535714860
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -535712410
$ws.Range("N122").Value = -13900
# Row 131
$ws.Range("H131").Value = 4940
$ws.Range("I131").Value = 995
$ws.Range("J131").Value = 6518
$ws.Range("K131").Value = 2985
$ws.Range("L131").Value = 19554
$ws.Range("M131").Value = 2055
$ws.Range("N131").Value = -29634
# Row 132
$ws.Range("H132").Value = 7147359.5
$ws.Range("I132").Value = 7940585
$ws.Range("J132").Value = 8329.429
$ws.Range("K132").Value = 23821755
$ws.Range("L132").Value = 24988.287
$ws.Range("M132").Value = -23819225
$ws.Range("N132").Value = -30048.287

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 16147196
$ws.Range("I32").Value = 18881248
$ws.Range("J32").Value = 46666.555
$ws.Range("K32").Value = 18881248
$ws.Range("L32").Value = 46666.555
$ws.Range("M32").Value = -18880961
$ws.Range("N32").Value = -47240.555
# Row 61
$ws.Range("H61").Value = 1828.711
$ws.Range("I61").Value = 1843.5
$ws.Range("J61").Value = 1792.3077
$ws.Range("K61").Value = 1843.5
$ws.Range("L61").Value = 1792.3077
$ws.Range("M61").Value = -1631.5
$ws.Range("N61").Value = -2216.3077
# Row 124
$ws.Range("H124").Value = 26143
$ws.Range("J124").Value = 26143
$ws.Range("L124").Value = 26143
$ws.Range("N124").Value = -35963
# Row 125
$ws.Range("H125").Value = 41715
$ws.Range("J125").Value = 41715
$ws.Range("L125").Value = 41715
$ws.Range("N125").Value = -51555
# Row 136
$ws.Range("H136").Value = 1828.711
$ws.Range("I136").Value = 1843.5
$ws.Range("J136").Value = 1792.3077
$ws.Range("K136").Value = 5530.5
$ws.Range("L136").Value = 5376.9231
$ws.Range("M136").Value = -2980.5
$ws.Range("N136").Value = -10476.9231

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2423.625
$ws.Range("I20").Value = 2365.1904
$ws.Range("J20").Value = 2535.182
$ws.Range("K20").Value = 2365.1904
$ws.Range("L20").Value = 2535.182
$ws.Range("M20").Value = -2118.1904
$ws.Range("N20").Value = -3029.182
# Row 86
$ws.Range("H86").Value = 1294678.5
$ws.Range("I86").Value = 2991
$ws.Range("J86").Value = 3324473.2
$ws.Range("K86").Value = 2991
$ws.Range("L86").Value = 3324473.2
$ws.Range("M86").Value = -1868
$ws.Range("N86").Value = -3326719.2
# Row 89
$ws.Range("H89").Value = 1294678.5
$ws.Range("I89").Value = 2991
$ws.Range("J89").Value = 3324473.2
$ws.Range("K89").Value = 14955
$ws.Range("L89").Value = 16622366
$ws.Range("M89").Value = -9339
$ws.Range("N89").Value = -16633598
# Row 105
$ws.Range("H105").Value = 22729206
$ws.Range("I105").Value = 1988.2727
$ws.Range("K105").Value = 1988.2727
$ws.Range("M105").Value = -241.2727
# Row 107
$ws.Range("H107").Value = 71429400
$ws.Range("I107").Value = 83334050
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 83334050
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = -83332130
$ws.Range("N107").Value = -5340
# Row 134
$ws.Range("H134").Value = 1526497.1
$ws.Range("I134").Value = 3171.291
$ws.Range("J134").Value = 6181104
$ws.Range("K134").Value = 9513.873
$ws.Range("L134").Value = 18543312
$ws.Range("M134").Value = -6978.873
$ws.Range("N134").Value = -18548382

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 129
$ws.Range("H129").Value = 15873973
$ws.Range("J129").Value = 33334674
$ws.Range("L129").Value = 100004022
$ws.Range("N129").Value = -100014022
# Row 130
$ws.Range("H130").Value = 62501468
$ws.Range("I130").Value = 166667310
$ws.Range("J130").Value = 1960
$ws.Range("K130").Value = 500001930
$ws.Range("L130").Value = 5880
$ws.Range("M130").Value = -499996910
$ws.Range("N130").Value = -15920
# Row 131
$ws.Range("H131").Value = 743.66
$ws.Range("J131").Value = 768.8570999999999
$ws.Range("L131").Value = 2306.5713
$ws.Range("N131").Value = -12386.5713

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4072.9062
$ws.Range("I70").Value = 4099.2104
$ws.Range("J70").Value = 4034.4614
$ws.Range("K70").Value = 4099.2104
$ws.Range("L70").Value = 4034.4614
$ws.Range("M70").Value = -3829.2104
$ws.Range("N70").Value = -4574.4614
# Row 73
$ws.Range("H73").Value = 4072.9062
$ws.Range("I73").Value = 4099.2104
$ws.Range("J73").Value = 4034.4614
$ws.Range("K73").Value = 4099.2104
$ws.Range("L73").Value = 4034.4614
$ws.Range("M73").Value = -3163.2104
$ws.Range("N73").Value = -5906.4614
# Row 102
$ws.Range("H102").Value = 1374.5385
$ws.Range("I102").Value = 1216.5
$ws.Range("J102").Value = 1510
$ws.Range("K102").Value = 1216.5
$ws.Range("L102").Value = 1510
$ws.Range("M102").Value = 405.5
$ws.Range("N102").Value = -4754
# Row 126
$ws.Range("H126").Value = 4225
$ws.Range("I126").Value = 3366.6667
$ws.Range("J126").Value = 5083.3335
$ws.Range("K126").Value = 10100.0001
$ws.Range("L126").Value = 15250.0005
$ws.Range("M126").Value = -7630.000100000001
$ws.Range("N126").Value = -20190.0005
# Row 132
$ws.Range("H132").Value = 6339.6
$ws.Range("I132").Value = 1638.6
$ws.Range("J132").Value = 15741.6
$ws.Range("K132").Value = 4915.799999999999
$ws.Range("L132").Value = 47224.8
$ws.Range("M132").Value = -2385.799999999999
$ws.Range("N132").Value = -52284.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 1593.5555
$ws.Range("I93").Value = 1606
$ws.Range("K93").Value = 1606
$ws.Range("M93").Value = -358

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 38462344
$ws.Range("I81").Value = 38462344
$ws.Range("K81").Value = 76924688
$ws.Range("M81").Value = -76923627
# Row 84
$ws.Range("H84").Value = 38462344
$ws.Range("I84").Value = 38462344
$ws.Range("K84").Value = 384623440
$ws.Range("M84").Value = -384618136
# Row 113
$ws.Range("H113").Value = 40000428
$ws.Range("I113").Value = 52631900
$ws.Range("K113").Value = 157895700
$ws.Range("M113").Value = -157893530
# Row 132
$ws.Range("H132").Value = 14943601
$ws.Range("I132").Value = 21762524
$ws.Range("J132").Value = 6913.7144
$ws.Range("K132").Value = 65287572
$ws.Range("L132").Value = 20741.1432
$ws.Range("M132").Value = -65285042
$ws.Range("N132").Value = -25801.1432
